# "No test guide version. Static v1.0 for all profiles"
#
# The "Mobile App Testing Guide Version" requirement row is retired (there is
# no longer a separate testing-guide-version question per generation), and
# the "Mobile App Specification Version" requirement now has a fixed answer
# of "v1.0" that applies to every profile instead of being left blank.

$wb = $excel.ActiveWorkbook
$gi = $wb.Worksheets.Item("General Information")
$files = $wb.Worksheets.Item("Files")

# Give the Mobile App Specification Version requirement (row 7) a static
# "v1.0" answer for every profile.
$gi.Range("D7").Value = "v1.0"

# Drop the now-unused "Mobile App Testing Guide Version" requirement row
# entirely (old row 8: Gen 7 / Yes / Mobile App Testing Guide Version) -
# everything below shifts up.
$gi.Rows(8).Delete()

# Keep the remaining generic requirement IDs sequential after the removal.
$gi.Range("A8").Value = "Gen 7"
$gi.Range("A9").Value = "Gen 8"

# Normalize the formatting on the rows touched above so they match the
# plain style used by the rest of the sheet (no stray reading-order
# override left behind from the old "Mobile App Testing Guide Version"
# row / the old trailing row).
$gi.Range("A2").Copy()
$gi.Range("A8:A9").PasteSpecial(-4122)
$gi.Range("B2").Copy()
$gi.Range("B7").PasteSpecial(-4122)
$gi.Range("C2").Copy()
$gi.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Files" sheet header row formatting normalized to match the rest of the
# header cells (no functional change to values).
$files.Range("A1").Copy()
$files.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
